$wb = $excel.ActiveWorkbook

# --- intInit sheet: rename a few rows and append new "...Gem"/"...Time" rows ---
$ws1 = $wb.Worksheets.Item("intInit")

$ws1.Range("A8").Value = "allianceFightTotalFightTime"
$ws1.Range("B8").Value = 60
$ws1.Rows.Item(8).RowHeight = 20

$ws1.Range("A9").Value = "allianceFightTimePerFight"
$ws1.Range("B9").Value = 5
$ws1.Rows.Item(9).RowHeight = 20

$ws1.Range("A7").Value = "allianceFightPrepareTime"
$ws1.Range("B7").Value = 30
$ws1.Rows.Item(7).RowHeight = 20

$ws1.Range("A6").Value = "activeShrineStageEventTime"
$ws1.Range("B6").Value = 15

$ws1.Range("A5").Value = "editAllianceTerrianHonour"
$ws1.Range("A4").Value = "editAllianceBasicInfoGem"
$ws1.Range("A3").Value = "buyArchonGem"
$ws1.Range("A2").Value = "createAllianceGem"

$ws1.Activate()
$ws1.Range("A3").Select()

# --- right sheet: selection moved off it, no longer the active tab ---
$ws3 = $wb.Worksheets.Item("right")
$ws3.Range("R3").Select()

$ws1.Activate()
